$p = $ppt.ActivePresentation
try {
  $w = $ppt.ActiveWindow
  Write-Host "window: $w"
} catch {
  Write-Host "no window: $_"
}
